$d = $word.ActiveDocument

# --- First paragraph: update paragraph formatting (border + indent) ---
$p1 = $d.Paragraphs(1)

# Add a paragraph border on all four sides with a 5pt gap to text, matching
# <w:pBdr><w:top w:space="5"/><w:left w:space="5"/><w:bottom w:space="5"/><w:right w:space="5"/></w:pBdr>
$p1.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromRight = 5

# Update the left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# --- First paragraph: update placeholder text & drop the trailing space run ---
# Replace "**ID__AFFARS_5342_topic_9__ID** " (including the trailing space
# that lived in its own run) with the new placeholder text, with no
# trailing space left behind.
$d.Content.Find.Execute("**ID__AFFARS_5342_topic_9__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_SUBPART_5342_15__ID**", 2)
